$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "BRENDAN BRYAN" heading -> split into "BRENDAN " + "BRYAN" runs,
# and drop the bold flags that were cached on the paragraph mark's rPr
# (the paragraph style's sz/szCs now also appear explicitly on the mark).
# ---------------------------------------------------------------------------
$namePara = $d.Paragraphs(1)
$nameXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="Name"/><w:jc w:val="center"/><w:rPr>' +
  '<w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/>' +
  '<w:sz w:val="34"/><w:szCs w:val="34"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/>' +
  '<w:b/><w:bCs/><w:sz w:val="34"/><w:szCs w:val="34"/></w:rPr>' +
  '<w:t xml:space="preserve">BRENDAN </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/>' +
  '<w:b/><w:bCs/><w:sz w:val="34"/><w:szCs w:val="34"/></w:rPr>' +
  '<w:t>BRYAN</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$null = $namePara.Range.InsertXML($nameXml)

# ---------------------------------------------------------------------------
# Change 2: the "Sharpened and Broadened Key Skills..." bullet is rewritten
# and split into two separate bullet paragraphs:
#   - "Practiced agile methodologies, "
#   - "Refreshed and improved " + "Key Skills and proficiencies with Udemy"
#     followed by the unchanged ", ChatGPT," and " and other online
#     resources." runs.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $par = $d.Paragraphs($i)
  if ($par.Range.Text -like "Sharpened and Broadened Key Skills*") {
    $target = $par
    break
  }
}

$skillsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="27"/></w:numPr>' +
  '<w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
  '<w:t xml:space="preserve">Practiced agile methodologies, </w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="27"/></w:numPr>' +
  '<w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
  '<w:t xml:space="preserve">Refreshed and improved </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
  '<w:t>Key Skills and proficiencies with Udemy</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
  '<w:t>, ChatGPT,</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
  '<w:t xml:space="preserve"> and other online resources.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

if ($target -ne $null) {
  $null = $target.Range.InsertXML($skillsXml)
}
